$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to remain text, since some values look like plain numbers
# and Excel would otherwise silently convert them to numeric values, losing the
# original text formatting (they were stored as inline strings in the workbook).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.554.72"
$ws.Range("D3").Value = "1.594.21"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "207.54"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "22.26"
$ws.Range("E8").Value = "  -4.04%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "1.822.30"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").Value = "1.610.13"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("E14").Value = "  -3.77%  "
$ws.Range("D15").Value = "0.538"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").Value = "63.32"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "27.540.74"
$ws.Range("D18").Value = "217.09"
$ws.Range("E18").Value = "  -4.84%  "
$ws.Range("D19").Value = "7.37"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").Value = "9.69"
$ws.Range("E23").Value = "  -3.67%  "
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "155.29"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "6.70"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("E29").Value = "  -4.49%  "
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("E32").Value = "  -3.57%  "
$ws.Range("D33").Value = "1.352.50"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("D34").Value = "2.95"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "0.957"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "0.0165"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("D40").Value = "0.814"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "0.954"
$ws.Range("E42").Value = "  -4.41%  "
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "63.91"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").Value = "1.731.88"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").Value = "87.27"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "0.0₇0992"
$ws.Range("E49").Value = "  -4.42%  "
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("E51").Value = "  -1.06%  "

# Restore the original (default) cell style now that the text values are set,
# so no residual number-format styling is left behind on column D.
$ws.Range("D2:D51").Style = "Normal"

